# Partially working snake game.
# Updates the Memory Layout Mapping tables on Sheet1:
#  - Hex range table (AC6:AE8): adjusts hex boundaries and adds a new
#    "SNAKE BODY" memory region row.
#  - Binary mirror table (AC16:AE18): follows the hex table via formulas.
#  - New decimal (FROM/TO) helper columns AG:AH that convert the binary
#    values back to decimal.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Hex "Memory Layout Mapping" table (AC6:AE8)
# ---------------------------------------------------------------

# Row 6: DATA + STACK now tops out at #af instead of #ef
$ws.Range("AD6").Value = "#af"

# Row 7: I/O PERIPHERAL range becomes #a0 - #af (was #f0 - #ff); re-style
# to match the plain centered look (style copied from AC15/From-To header).
$ws.Range("AC15").Copy() | Out-Null
$ws.Range("AC7:AD7").PasteSpecial(-4122) | Out-Null
$ws.Range("AC7").Value = "#a0"
$ws.Range("AD7").Value = "#af"

# Row 8: brand new SNAKE BODY region, #c0 - #ff.
$ws.Range("AC6").Copy() | Out-Null
$ws.Range("AC8:AD8").PasteSpecial(-4122) | Out-Null
$ws.Range("AC8").Value = "#c0"
$ws.Range("AD8").Value = "#ff"
$ws.Range("AC15").Copy() | Out-Null
$ws.Range("AE8").PasteSpecial(-4122) | Out-Null
$ws.Range("AE8").Value = "SNAKE BODY"

# ---------------------------------------------------------------
# Binary mirror table (AC16:AE18) - driven by DEC2BIN/HEX2DEC formulas
# ---------------------------------------------------------------

# AC17:AD18 fill down/across as one shared formula block referencing the
# (now 3-row) hex table above. Row 18 starts blank, so copy the binary
# cell format down first, then overwrite with the shared formula.
$ws.Range("AC17:AD17").Copy() | Out-Null
$ws.Range("AC18:AD18").PasteSpecial(-4122) | Out-Null
$ws.Range("AC17:AD18").Formula = "=DEC2BIN(HEX2DEC(MID(AC7, 2, LEN(AC7))))"

# New SNAKE BODY row in the binary table.
$ws.Range("AC15").Copy() | Out-Null
$ws.Range("AE18").PasteSpecial(-4122) | Out-Null
$ws.Range("AE18").Value = "SNAKE BODY"

# ---------------------------------------------------------------
# New decimal helper columns (AG:AH) - convert each binary string back
# to its decimal value for double-checking.
# ---------------------------------------------------------------

$ws.Range("AC15").Copy() | Out-Null
$ws.Range("AG15:AH15").PasteSpecial(-4122) | Out-Null
$ws.Range("AG15").Value = "FROM"
$ws.Range("AH15").Value = "TO"

$ws.Range("AG16:AG18").Formula = "=SUMPRODUCT(MID(AC16,LEN(AC16)-ROW(INDIRECT(""1:""&LEN(AC16)))+1,1)*2^(ROW(INDIRECT(""1:""&LEN(AC16)))-1))"

$ws.Range("AH16").Formula = "=SUMPRODUCT(MID(AD16,LEN(AD16)-ROW(INDIRECT(""1:""&LEN(AD16)))+1,1)*2^(ROW(INDIRECT(""1:""&LEN(AD16)))-1))"
$ws.Range("AH17:AH18").Formula = "=SUMPRODUCT(MID(AD17,LEN(AD17)-ROW(INDIRECT(""1:""&LEN(AD17)))+1,1)*2^(ROW(INDIRECT(""1:""&LEN(AD17)))-1))"

# ---------------------------------------------------------------
# View bookkeeping to match the workbook as last saved by the author.
# ---------------------------------------------------------------
$ws.Range("AE12").Select() | Out-Null

$excel.CalculateFull() | Out-Null
